$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text (shared-string) cell values in the exact order needed so that the
# generated sharedStrings.xml table indexes match the target workbook.
$ws.Range("C6").Value = " We\'re considering going after\na high-ranking outlaw. We just may go for it."
$ws.Range("A6").Value = "SCRIPT/G01P03A/um1111.ssb"
$ws.Range("C7").Value = " It would be a tough fight. But\nthere\'s only one enemy."
$ws.Range("C8").Value = " With all three of us in on it,\nwe can take one outlaw!"
$ws.Range("D6").Value = " Мы собираемся на охоту за\nвысокоранговым негодяем."
$ws.Range("D7").Value = " Это будет нелегкий бой. Но\nвраг только один."
$ws.Range("D8").Value = " Если мы все трое сразимся с\nним, мы сможем справиться с одним\nнегодяем!"
$ws.Range("E6").Value = " Íú òïáéñàåíòÿ îà ïöïóô èà\nâúòïëïñàîãïâúí îåãïäÿåí."
$ws.Range("E7").Value = " Üóï áôäåó îåìåãëéê áïê. Îï\nâñàã óïìûëï ïäéî."
$ws.Range("E8").Value = " Åòìé íú âòå óñïå òñàèéíòÿ ò\nîéí, íú òíïçåí òðñàâéóûòÿ ò ïäîéí\nîåãïäÿåí!"
$ws.Range("A7").Value = "SCRIPT/G01P03A/um1116.ssb"

# Numeric (non shared-string) cell values
$ws.Range("B6").Value = 313
$ws.Range("B7").Value = 316
$ws.Range("B8").Value = 319

# --- Row 8 style: new thin-bottom border + smaller (8pt) font like the rows above ---
$row8 = $ws.Range("A8:E8")
$row8.WrapText = $true
$row8.Borders.Item(9).LineStyle = 1
$ws.Range("C8:E8").Font.Size = 8

# --- Row heights ---
$ws.Rows.Item(6).RowHeight = 43.2
$ws.Rows.Item(7).RowHeight = 43.2
$ws.Rows.Item(8).RowHeight = 21.6

# --- Update selection to D6 ---
$ws.Range("D6").Select() | Out-Null
